$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 185-186; this pushes the existing rows
# 185-281 down to 187-283 (and Excel auto-extends the dimension).
$ws.Rows("185:186").Insert()

# Row 185 - "Primera" quality, new week's data
$ws.Range("A185").Value = 6
$ws.Range("B185").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C185").Value = "Metropolitana"
$ws.Range("D185").Value = 44455
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 100112039
$ws.Range("G185").Value = "Ciboulette"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 750
$ws.Range("K185").Value = 1000
$ws.Range("L185").Value = 1200
$ws.Range("M185").Value = 1117
$ws.Range("N185").Value = "$/docena de atados"
$ws.Range("O185").Value = "Región Metropolitana"
$ws.Range("P185").Value = 372
$ws.Range("Q185").Value = 3
$ws.Range("R185").Value = "Hortaliza"

# Row 186 - "Segunda" quality, new week's data
$ws.Range("A186").Value = 6
$ws.Range("B186").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C186").Value = "Metropolitana"
$ws.Range("D186").Value = 44455
$ws.Range("E186").Value = 13
$ws.Range("F186").Value = 100112039
$ws.Range("G186").Value = "Ciboulette"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Segunda"
$ws.Range("J186").Value = 250
$ws.Range("K186").Value = 800
$ws.Range("L186").Value = 800
$ws.Range("M186").Value = 800
$ws.Range("N186").Value = "$/docena de atados"
$ws.Range("O186").Value = "Región Metropolitana"
$ws.Range("P186").Value = 267
$ws.Range("Q186").Value = 3
$ws.Range("R186").Value = "Hortaliza"
